$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.458.98'
$ws.Range("E2").Value = '  -1.13%  '

$ws.Range("D3").Value = '2.337.97'
$ws.Range("E3").Value = '  +3.02%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.649'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.14'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.62%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.455'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.50%  '

$ws.Range("E10").Value = '  -3.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.33%  '

$ws.Range("D13").Value = '2.684.32'
$ws.Range("E13").Value = '  +2.92%  '

$ws.Range("E14").Value = '  -1.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.853'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.64%  '

$ws.Range("D18").Value = '2.328.43'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("D19").Value = '43.378.70'
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("E20").Value = '  -2.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.09%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("E27").Value = '  -1.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.81%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.129'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.14%  '

$ws.Range("E33").Value = '  +0.65%  '

$ws.Range("E34").Value = '  +4.55%  '

$ws.Range("E35").Value = '  -1.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.51'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.31%  '

$ws.Range("E40").Value = '  -1.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.83%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.62%  '

$ws.Range("E44").Value = '  +8.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("E46").Value = '  +0.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0948'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("D49").Value = '1.447.62'
$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("B50").Value = 'TerraClassic'
$ws.Range("C50").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000205'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.60%  '

$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.90%  '
